# "Se crea la lista de alumnos" -- reset the roster sheet back to a blank
# template: wipe out the student rows (columns B:E, rows 2-10) that held
# names / aulas / mails / github handles, drop the mailto: hyperlinks that
# lived on D2:D10, and leave the selection on D10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the mailto hyperlinks attached to the "mail" column before wiping
# the cell text itself (ClearContents alone would leave the <hyperlinks>
# entries dangling).
$ws.Range("D2:D10").Hyperlinks.Delete()

# Blank out the student data (Nombre y Apellido, Aula, mail, nombre github)
# for the 9 existing rows while keeping the row/column formatting intact.
$ws.Range("B2:E10").ClearContents()

# Leave the selection where the author left it.
[void]$ws.Range("D10").Select()
